$wb = $excel.ActiveWorkbook

# --- Step 1: Reorder sheets ---
# Current order: Sheet, Web Development 1, Cryptography 1, Python 1, General 2, General Easy 1, General 3
# Target order (before adding Python 3): Sheet, Web Development 1, Python 1, General 2, General Easy 1, General 3, Cryptography 1
$cryptoWs = $wb.Worksheets.Item("Cryptography 1")
$generalThreeWs = $wb.Worksheets.Item("General 3")
$cryptoWs.Move($null, $generalThreeWs)

# --- Step 2: Trim "General 3" - drop the first two quiz questions, shift rows up ---
$g3 = $wb.Worksheets.Item("General 3")
$g3.Rows.Item(2).Delete()
$g3.Rows.Item(2).Delete()
$g3.Cells.Item(1,2).Value = 'A directory'
$g3.Cells.Item(1,3).Value = 'An access code'

# --- Step 3: Trim "Cryptography 1" - drop the first two quiz questions, shift rows up ---
$crypto = $wb.Worksheets.Item("Cryptography 1")
$crypto.Rows.Item(2).Delete()
$crypto.Rows.Item(2).Delete()
$crypto.Cells.Item(1,2).Value = 'Translation'
$crypto.Cells.Item(1,3).Value = 'Decryption'

# --- Step 4: Add new "Python 3" sheet at the end, populate with quiz data ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newWs = $wb.Worksheets.Add($null, $lastSheet)
$newWs.Name = 'Python 3'
$newWs.Cells.Item(1,1).Value = 'Python 3'
$newWs.Cells.Item(1,2).Value = 'Hard'
$newWs.Cells.Item(1,3).Value = 'Python'
$newWs.Cells.Item(2,1).Value = 'What does the @classmethod decorator do?'
$newWs.Cells.Item(2,2).Value = 'Converts a method into a static method'
$newWs.Cells.Item(2,3).Value = 'Converts a method so it can only be called on class instances'
$newWs.Cells.Item(2,4).Value = 'Converts a method to receive the class as the first argument instead of the instance'
$newWs.Cells.Item(2,5).Value = 'Makes the method private'
$newWs.Cells.Item(2,6).Value = 2
$newWs.Cells.Item(3,1).Value = 'Which of the following statements is TRUE about Python’s GIL (Global Interpreter Lock)?'
$newWs.Cells.Item(3,2).Value = 'It allows true parallel execution of threads on multiple CPUs'
$newWs.Cells.Item(3,3).Value = 'It prevents memory leaks in C extensions'
$newWs.Cells.Item(3,4).Value = 'It ensures only one thread executes Python bytecode at a time'
$newWs.Cells.Item(3,5).Value = 'It is removed in Python 3.10'
$newWs.Cells.Item(3,6).Value = 2
$newWs.Cells.Item(4,1).Value = 'How can you create a custom context manager in Python?'
$newWs.Cells.Item(4,2).Value = 'Using @classmethod'
$newWs.Cells.Item(4,3).Value = 'Using __init__ and __del__'
$newWs.Cells.Item(4,4).Value = 'By defining __enter__ and __exit__ methods'
$newWs.Cells.Item(4,5).Value = 'By using yield without any decorators'
$newWs.Cells.Item(4,6).Value = 2
$newWs.Cells.Item(5,1).Value = 'Which of the following is used to define an abstract method in Python?'
$newWs.Cells.Item(5,2).Value = 'from abc import abstractmethod'
$newWs.Cells.Item(5,3).Value = '@abstractmethod'
$newWs.Cells.Item(5,4).Value = 'def method(self): pass'
$newWs.Cells.Item(5,5).Value = 'Both A and B'
$newWs.Cells.Item(5,6).Value = 3
$newWs.Cells.Item(6,1).Value = 'What is the primary purpose of __slots__ in a class?'
$newWs.Cells.Item(6,2).Value = 'To define class-level constants'
$newWs.Cells.Item(6,3).Value = 'To create a dictionary for instance variables'
$newWs.Cells.Item(6,4).Value = 'To prevent dynamic creation of new attributes and reduce memory'
$newWs.Cells.Item(6,5).Value = 'To define abstract methods'
$newWs.Cells.Item(6,6).Value = 2
$newWs.Cells.Item(7,1).Value = 'Which statement is true about Python descriptors?'
$newWs.Cells.Item(7,2).Value = 'They only work with class methods'
$newWs.Cells.Item(7,3).Value = 'They define attribute access using __get__, __set__, __delete__'
$newWs.Cells.Item(7,4).Value = 'Descriptors are only available in Python 3.10+'
$newWs.Cells.Item(7,5).Value = 'Descriptors cannot be reused across classes'
$newWs.Cells.Item(7,6).Value = 1
$newWs.Cells.Item(8,1).Value = 'What is a coroutine in Python?'
$newWs.Cells.Item(8,2).Value = 'A thread'
$newWs.Cells.Item(8,3).Value = 'A class with __iter__'
$newWs.Cells.Item(8,4).Value = 'a function paused with yield and resumed later'
$newWs.Cells.Item(8,5).Value = 'A method in a metaclass'
$newWs.Cells.Item(8,6).Value = 2
$newWs.Cells.Item(9,1).Value = 'Which of the following is not a valid use of a metaclass?'
$newWs.Cells.Item(9,2).Value = 'Injecting methods into a class'
$newWs.Cells.Item(9,3).Value = 'Enforcing naming conventions'
$newWs.Cells.Item(9,4).Value = 'Changing inheritance at runtime'
$newWs.Cells.Item(9,5).Value = 'Dynamically creating instance variables'
$newWs.Cells.Item(9,6).Value = 3
$newWs.Cells.Item(10,1).Value = 'What does weakref module provide in Python?'
$newWs.Cells.Item(10,2).Value = 'Multithreading support'
$newWs.Cells.Item(10,3).Value = 'References that do not increase reference count'
$newWs.Cells.Item(10,4).Value = 'Dynamic type checking'
$newWs.Cells.Item(10,5).Value = 'Object pooling'
$newWs.Cells.Item(10,6).Value = 1
$newWs.Cells.Item(11,1).Value = 'What does the super() function do in Python?'
$newWs.Cells.Item(11,2).Value = 'Calls a method from a subclass'
$newWs.Cells.Item(11,3).Value = 'Calls the next method in the method resolution order (MRO)'
$newWs.Cells.Item(11,4).Value = 'Creates a static method'
$newWs.Cells.Item(11,5).Value = 'Returns the parent class constructor  directly'
$newWs.Cells.Item(11,6).Value = 1

